# "Generate Report for Archive" — refresh the localization-status report:
# the handoff that was pending is now in translation, so update the status
# cells on all three sheets and re-fit the columns that held the
# (now shorter) status text, same as Excel does automatically when the
# cell content that drove the column's width changes.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: per-language status columns E (zh-cn) and F (de-de)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Per-language detail sheets: Status column C
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Re-fit the columns whose text just changed length (they shrink now that
# "In Translation" is shorter than "Ready for handoff").
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null

# The host's AutoFit only coarsely approximates real Excel's font-metric
# based pixel fit, so nail the resulting column width down explicitly to
# match what Excel actually computed for the new, shorter status text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
